$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Wed Nov 01 15:53:13 EDT 2023"
$ws.Range("B3").Value = "Wed Nov 01 15:53:25 EDT 2023"
$ws.Range("B4").Value = "Wed Nov 01 15:53:36 EDT 2023"
$ws.Range("B5").Value = "Wed Nov 01 15:53:48 EDT 2023"
